# Update Frame class and documentation.
#
# Adds a new "Frame:" subsection (heading + description paragraph) in both
# the high-level "Class Design" overview and the detailed "Method Design"
# section, and marks the rendered-page-break hint that Word recorded after
# the new content pushed the "Traffic Light Rectangle" method description
# onto a new page.

$d = $word.ActiveDocument

# Pull the whole package as WordprocessingML so we can make precise,
# surgical edits to the underlying markup and feed it back in one shot.
$full = $d.Content.XML(1)

# --- Change 1: overview "Class Design" section --------------------------
# The empty paragraph that used to sit between the "Traffic Light
# Rectangle:" description and the "Method Design" heading becomes the new
# "Frame:" heading, followed by its description paragraph, followed by a
# fresh blank paragraph (taking over the old spacer's role).
$old1 = '<w:p w14:paraId="0CA2525F" w14:textId="77777777" w:rsidR="00C56331" w:rsidRPr="00C56331" w:rsidRDefault="00C56331" w:rsidP="00C56331"><w:pPr><w:ind w:left="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'

$new1 = '<w:p w14:paraId="0CA2525F" w14:textId="77777777" w:rsidR="00C56331" w:rsidRPr="00C56331" w:rsidRDefault="00C56331" w:rsidP="00C56331">' `
  + '<w:pPr><w:ind w:left="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>Frame:</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p><w:pPr><w:ind w:left="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The frame class will act as the frame window for the GUI. </w:t></w:r>' `
  + '</w:p>' `
  + '<w:p><w:pPr><w:ind w:left="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'

if (([regex]::Matches($full, [regex]::Escape($old1))).Count -ne 1) {
  throw "anchor 1 not uniquely found"
}
$full = $full.Replace($old1, $new1)

# --- Change 2: rendered-page-break hint ----------------------------------
# The new content above pushes the "Traffic Light Rectangle" method
# description (detailed "Method Design" section) onto a fresh page.
$old2 = '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The class contains an override method to edit graphic components of the class. </w:t></w:r></w:p><w:p w14:paraId="093AFE9E"'

$new2 = '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">The class contains an override method to edit graphic components of the class. </w:t></w:r></w:p><w:p w14:paraId="093AFE9E"'

if (([regex]::Matches($full, [regex]::Escape($old2))).Count -ne 1) {
  throw "anchor 2 not uniquely found"
}
$full = $full.Replace($old2, $new2)

# --- Change 3: detailed "Method Design" section --------------------------
# The final paragraph (which carries the _GoBack bookmark) becomes the new
# "Frame:" heading; the description paragraph follows keeping the bookmark
# at its end, then a fresh trailing blank paragraph closes the document.
$old3 = '<w:p w14:paraId="093AFE9E" w14:textId="77777777" w:rsidR="0073029F" w:rsidRPr="0073029F" w:rsidRDefault="0073029F" w:rsidP="0033745E"><w:pPr><w:ind w:left="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$new3 = '<w:p w14:paraId="093AFE9E" w14:textId="77777777" w:rsidR="0073029F" w:rsidRPr="0073029F" w:rsidRDefault="0073029F" w:rsidP="0033745E">' `
  + '<w:pPr><w:ind w:left="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>Frame:</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p><w:pPr><w:ind w:left="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The Frame class contains methods for actions to be performed such as dialogue box popups and exiting the window. </w:t></w:r>' `
  + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' `
  + '</w:p>' `
  + '<w:p><w:pPr><w:ind w:left="720"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'

if (([regex]::Matches($full, [regex]::Escape($old3))).Count -ne 1) {
  throw "anchor 3 not uniquely found"
}
$full = $full.Replace($old3, $new3)

# Write the fully patched package back into the document in one shot.
$d.Content.InsertXML($full)
